$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-52:
# the serial date value moves from 45204 (2023-10-05) to 45205 (2023-10-06).
for ($row = 2; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
